$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.777.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.234.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.648"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "229.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.52%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.444"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0954"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "56.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.38"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.27%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.566.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.821"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.224.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.652.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0978"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.20%  "
$ws.Range("E21").Value = "  -1.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.79%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -6.22%  "
$ws.Range("E26").Value = "  -3.04%  "
$ws.Range("E27").Value = "  +23.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0695"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.74"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.90%  "
$ws.Range("E38").Value = "  -2.81%  "
$ws.Range("E39").Value = "  -5.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0257"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.000218"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "96.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("E47").Value = "  -2.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  +1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.427.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.69%  "
$ws.Range("E51").Value = "  +1.72%  "
